$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 29 mirrors the existing "request" rows (e.g. row 28): a new
# item request submitted by the same requester (Shahzeb Ahmed), this
# time for the OPP committee's homecoming event, item A12 (Aluminum),
# quantity 3.
$ws.Range("B29").Value = "shahzeb2000786@gmail.com"
$ws.Range("C29").Value = "Shahzeb Ahmed"
$ws.Range("D29").Value = "OPP"
$ws.Range("E29").Value = "homecoming"
$ws.Range("F29").Value = "A12"
$ws.Range("G29").Value = "Aluminum"

# H29 ("Quantity") is stored as text (like every other cell in this
# column), not a number -- format as Text first so Excel doesn't coerce
# "3" into a numeric cell, then drop the format back to the sheet's
# default "Normal" style so no stray cell-style survives the edit.
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "3"
$ws.Range("H29").Style = "Normal"
